# Update countries & provincias Spain
# Reflects the COVID data refresh dated 3 de Octubre de 2020 a las 19:49

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 3 de Octubre de 2020 a las 19:49"

# --- Update Estados Unidos (row 4) ---
$ws.Range("B4").Value = 7570742
$ws.Range("C4").Value = 21419
$ws.Range("D4").Value = 4782593
$ws.Range("E4").Value = 2574391
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 234
$ws.Range("H4").Value = 213758

# --- Update India (row 5) ---
$ws.Range("B5").Value = 6543448
$ws.Range("C5").Value = 71514
$ws.Range("D5").Value = 5501813
$ws.Range("E5").Value = 939913
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 847
$ws.Range("H5").Value = 101722

# --- Turquia overtakes Italia: rows 21/22 swap order ---
# Row 21 becomes Turquia (with its refreshed totals)
$ws.Range("A21").Value = "Turquia"
$ws.Range("B21").Value = 323014
$ws.Range("C21").Value = 1502
$ws.Range("D21").Value = 283868
$ws.Range("E21").Value = 30762
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 59
$ws.Range("H21").Value = 8384

# Row 22 becomes Italia (keeping its previous totals)
$ws.Range("A22").Value = "Italia"
$ws.Range("B22").Value = 322751
$ws.Range("C22").Value = 2844
$ws.Range("D22").Value = 231217
$ws.Range("E22").Value = 55566
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 27
$ws.Range("H22").Value = 35968

# --- Update Irlanda (row 74) ---
$ws.Range("B74").Value = 37668
$ws.Range("C74").Value = 605
$ws.Range("D74").Value = 23364
$ws.Range("E74").Value = 12494
$ws.Range("F74").Value = 0
$ws.Range("G74").Value = 9
$ws.Range("H74").Value = 1810

# --- Update Sri Lanka (row 143) ---
$ws.Range("B143").Value = 3395
$ws.Range("C143").Value = 7
$ws.Range("D143").Value = 3254
$ws.Range("E143").Value = 128
$ws.Range("F143").Value = 0
$ws.Range("G143").Value = 0
$ws.Range("H143").Value = 13

# --- Update Mali (row 146) ---
$ws.Range("B146").Value = 3170
$ws.Range("C146").Value = 14
$ws.Range("D146").Value = 2476
$ws.Range("E146").Value = 563
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 131

# --- Update Islas Feroe (row 179) ---
$ws.Range("B179").Value = 473
$ws.Range("C179").Value = 1
$ws.Range("D179").Value = 432
$ws.Range("E179").Value = 41
$ws.Range("F179").Value = 0
$ws.Range("G179").Value = 0
$ws.Range("H179").Value = 0

# --- Santa Lucia / Nueva Caledonia swap order (tied values, rows 207/208) ---
$ws.Range("A207").Value = "Santa Lucia"
$ws.Range("B207").Value = 27
$ws.Range("C207").Value = 0
$ws.Range("D207").Value = 27
$ws.Range("E207").Value = 0
$ws.Range("F207").Value = 0
$ws.Range("G207").Value = 0
$ws.Range("H207").Value = 0

$ws.Range("A208").Value = "Nueva Caledonia"
$ws.Range("B208").Value = 27
$ws.Range("C208").Value = 0
$ws.Range("D208").Value = 27
$ws.Range("E208").Value = 0
$ws.Range("F208").Value = 0
$ws.Range("G208").Value = 0
$ws.Range("H208").Value = 0
